$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of trade data (row 15) - retrieved via the new fundamentals method
$ws.Range("A15").Value = 10627.62
$ws.Range("B15").Value = 10413.11
$ws.Range("C15").Value = 104.49
$ws.Range("D15").Value = 106.64
$ws.Range("E15").Value = $false
$ws.Range("F15").Value = 2.06
$ws.Range("G15").Value = 42626.544479166667
$ws.Range("G15").NumberFormat = "m/d/yy h:mm"
$ws.Range("H15").Value = $true
